$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns I ("I0") and J ("IF") — header cells copy H1's style (bold,
# bordered, centered) so they match the existing header formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data: row, I value, J value
$data = @(
    @(2, 8, 8),
    @(3, 8, 8),
    @(4, 6, 6),
    @(5, 7, 8),
    @(6, 6, 6),
    @(7, 9, 9),
    @(8, 7, 7),
    @(9, 6, 6),
    @(10, 8, 8),
    @(11, 7, 8),
    @(12, 7, 7),
    @(13, 7, 7),
    @(14, 9, 9),
    @(15, 6, 6),
    @(16, 8, 8),
    @(17, 5, 6),
    @(18, 6, 7),
    @(19, 5, 5),
    @(20, 10, 10),
    @(21, 7, 7),
    @(22, 6, 6),
    @(23, 9, 9),
    @(24, 8, 8),
    @(25, 5, 6),
    @(26, 8, 8),
    @(27, 7, 7),
    @(28, 6, 6),
    @(29, 6, 6),
    @(30, 6, 6),
    @(31, 5, 5),
    @(32, 7, 7),
    @(33, 9, 9),
    @(34, 8, 8),
    @(35, 8, 8),
    @(36, 8, 8),
    @(37, 7, 8),
    @(38, 8, 8),
    @(39, 7, 8),
    @(40, 8, 9),
    @(41, 7, 7),
    @(42, 5, 6),
    @(43, 5, 6),
    @(44, 7, 8),
    @(45, 8, 8),
    @(46, 10, 10),
    @(47, 7, 7),
    @(48, 5, 7),
    @(49, 8, 9),
    @(50, 7, 8),
    @(51, 7, 8),
    @(52, 7, 7),
    @(53, 8, 8),
    @(54, 6, 7),
    @(55, 6, 6),
    @(56, 8, 8),
    @(57, 8, 8),
    @(58, 7, 7),
    @(59, 6, 7),
    @(60, 6, 6),
    @(61, 1, 2),
    @(62, 7, 8),
    @(63, 3, 4),
    @(64, 7, 7),
    @(65, 6, 6)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}
